# Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
# Swap the content of columns B:AD between paired rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(61, 62),
    @(156, 157),
    @(228, 229),
    @(305, 306)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AD$r1")
    $range2 = $ws.Range("B$r2" + ":AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
